# Weekly update: prepend a new week's worth of data (2 rows) to the
# "Perejil" (parsley) price table, pushing the existing rows down by two.
#
# The block of daily/weekly observations in this sheet lives in rows
# 2-51 (header in row 1). A new pair of rows is inserted right before the
# current row 30, and the previously-existing rows 30-51 shift down to
# 32-53. The new rows 30 and 31 get this week's figures; all of the
# "constant" columns (market, region, category, quality, origin,
# classification, etc.) are identical throughout this block, so they are
# copied from the row that used to be row 30 (now row 32) and then the
# handful of cells that actually change are overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new blank rows above the current row 30. This shifts the
#    old rows 30-51 down to 32-53, preserving their contents/styles.
$ws.Range("A30:R31").EntireRow.Insert()

# 2) Seed the two new rows with the same "constant" values as the rest of
#    the block (taken from the row that is now at 32, i.e. the old row 30).
$ws.Range("A30:R30").Value = $ws.Range("A32:R32").Value()
$ws.Range("A31:R31").Value = $ws.Range("A32:R32").Value()

# 3) Overwrite the cells that hold this week's new observations.

# Row 30 (new)
$ws.Range("D30").Value = 45161
$ws.Range("J30").Value = 285
$ws.Range("M30").Value = 1905
$ws.Range("N30").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("P30").Value = 1270
$ws.Range("Q30").Value = 1.5

# Row 31 (new)
$ws.Range("D31").Value = 45161
$ws.Range("J31").Value = 420
$ws.Range("K31").Value = 1800
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = 1895
$ws.Range("P31").Value = 948
